# Apply "Added handling of common packages." change:
# Reorders the Field Name / Field Type rows for each class on the
# "classFields" sheet so that fields belonging to a common/shared
# package layout line up consistently across classes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

# Final (Field Name, Field Type) pairs for rows 2-17 (Class Name and
# Field Modifier columns - A and C - are unchanged).
$rows = @{
    2  = @("availableItems", "int")
    3  = @("reservedItems", "int")
    4  = @("name", "java.lang.String")
    5  = @("id", "java.lang.Long")
    6  = @("repository", "com.zatribune.spring.ecommerce.stock.db.repository.ProductRepository")
    7  = @("log", "org.slf4j.Logger")
    8  = @("reservedItems", "int")
    9  = @("id", "java.lang.Long")
    10 = @("availableItems", "int")
    11 = @("name", "java.lang.String")
    12 = @("log", "org.slf4j.Logger")
    13 = @("orderService", "com.zatribune.spring.ecommerce.stock.service.OrderService")
    14 = @("repository", "com.zatribune.spring.ecommerce.stock.db.repository.ProductRepository")
    15 = @("SOURCE", "domain.OrderSource")
    16 = @("template", "org.springframework.kafka.core.KafkaTemplate")
    17 = @("log", "org.slf4j.Logger")
}

foreach ($r in $rows.Keys) {
    $pair = $rows[$r]
    $ws.Cells.Item($r, 2).Value = $pair[0]
    $ws.Cells.Item($r, 4).Value = $pair[1]
}
